$d = $word.ActiveDocument

$replacements = @(
    @("81÷5=", "84÷9="),
    @("86÷3=", "47÷7="),
    @("46÷3=", "34÷7="),
    @("76÷8=", "63÷4="),
    @("98÷2=", "28÷4="),
    @("35÷3=", "60÷2="),
    @("42÷6=", "41÷4="),
    @("42÷7=", "73÷7="),
    @("92÷9=", "49÷5="),
    @("72÷9=", "21÷8="),
    @("77÷7=", "12÷8="),
    @("61÷9=", "62÷3="),
    @("35÷6=", "80÷4="),
    @("28÷7=", "39÷7="),
    @("33÷7=", "28÷3="),
    @("11÷5=", "95÷9="),
    @("27÷5=", "91÷7="),
    @("54÷6=", "49÷4="),
    @("64÷5=", "80÷8="),
    @("42÷4=", "14÷6="),
    @("57÷9=", "85÷8="),
    @("52÷8=", "81÷4="),
    @("52÷7=", "83÷5="),
    @("85÷2=", "39÷6="),
    @("52÷9=", "97÷2=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
